# Adapt column header formatting to respective input file names (#7)
#
# Renames the "_old" / "_new" header-name suffixes to "_FV2210" / "_FV2304"
# (matching the new naming scheme keyed off the EDIFACT format version),
# wraps the sheet's data range in a native Excel Table so the new headers
# are also reflected as table column names, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells in row 1 (A1:J1 = "_old" -> "_FV2210", L1:U1 = "_new" -> "_FV2304") ---
$oldSuffixCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newSuffixCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value() -replace "_old$", "_FV2210")
}

foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value() -replace "_new$", "_FV2304")
}

# --- 2. Turn the used range into a proper Excel Table (ListObject) ---
$tableRange = $ws.Range("A1:U85")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (pane split below row 1) ---
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
